# Insert a new weekly price record as row 56 on the (single) sheet,
# pushing all existing rows 56-171 down to 57-172.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value  = 4
$ws.Cells.Item(56, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(56, 3).Value  = 'Los Lagos'
$ws.Cells.Item(56, 4).Value  = 44536
$ws.Cells.Item(56, 5).Value  = 10
$ws.Cells.Item(56, 6).Value  = 'Fruta'
$ws.Cells.Item(56, 7).Value  = 100108
$ws.Cells.Item(56, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(56, 9).Value  = 100108005
$ws.Cells.Item(56, 10).Value = 'Piña'
$ws.Cells.Item(56, 11).Value = 'Caramelo'
$ws.Cells.Item(56, 12).Value = 'Tercera'
$ws.Cells.Item(56, 13).Value = 120
$ws.Cells.Item(56, 14).Value = 20000
$ws.Cells.Item(56, 15).Value = 21000
$ws.Cells.Item(56, 16).Value = 20500
$ws.Cells.Item(56, 17).Value = '$/caja 16 unidades'
$ws.Cells.Item(56, 18).Value = 'Ecuador'
$ws.Cells.Item(56, 19).Value = 1281
$ws.Cells.Item(56, 20).Value = 16
